$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.542.91"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.825.16"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.05%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.74"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.51%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +0.06%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5108"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -5.50%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3951"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -1.34%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08236"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +6.23%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.114"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -0.59%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.80"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -0.50%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.361"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +0.23%  "
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.17"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  +0.01%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.566"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "1.818.88"
$ws.Range("E16").Value = "  -0.55%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001124"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +3.11%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.00"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +3.42%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06662"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +1.14%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.83"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  +0.03%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.107"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "28.584.70"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +2.09%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.272"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +0.36%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.36"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +2.54%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.11"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").Value = "2.029.72"
$ws.Range("E28").Value = "  -0.47%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.412"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -1.71%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.73"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +1.82%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.118"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -1.56%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1088"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -2.96%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.787"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("E34").Value = "  +0.39%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07075"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -4.54%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2233"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -1.19%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02357"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -0.15%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.254"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +0.69%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.778"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -2.05%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6357"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +0.88%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.28"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -1.21%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.183"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -0.70%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.402"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -0.21%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.57"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +0.03%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5977"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +1.20%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.739"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("E47").Value = "  -0.07%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.994"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("E49").Value = "  -0.12%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06945"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +0.36%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.086"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +4.35%  "
